$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览): simple numeric updates ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range('F6').Value = 3199
$ws1.Range('F7').Value = 2768
$ws1.Range('F11').Value = 359
$ws1.Range('F12').Value = 299
$ws1.Range('F13').Value = 31
$ws1.Range('F14').Value = 5790
$ws1.Range('F18').Value = 167
$ws1.Range('F20').Value = 480
$ws1.Range('F22').Value = 76
$ws1.Range('F24').Value = 2018
$ws1.Range('F26').Value = 338

# --- Sheet 2 (演出): simple numeric updates ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range('F25').Value = 4039
$ws2.Range('F33').Value = 3

# --- Sheet 3 (本地生活): simple numeric updates ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range('F5').Value = 2574
$ws3.Range('F13').Value = 568

# --- Sheet 4 (全部类型): simple numeric updates ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range('F14').Value = 3199
$ws4.Range('F15').Value = 2768
$ws4.Range('F19').Value = 568
$ws4.Range('F20').Value = 359
$ws4.Range('F23').Value = 299
$ws4.Range('F24').Value = 31
$ws4.Range('F44').Value = 2018
$ws4.Range('F48').Value = 338

# --- Sheet 4: rows 25, 28-42 content updates (regenerated aggregate rows) ---
# Row 25
$ws4.Range('B25').Value = '2024.03.09'
$ws4.Range('C25').Value = '上海·S·CGE动漫游戏嘉年华'
$ws4.Range('D25').Value = '军工路1076号 纪希片场(秀场)'
$ws4.Range('E25').Value = '2024.03.09 10:00-03.10 17:00'
$ws4.Range('F25').Value = 5790
$ws4.Range('G25').Value = 70
$ws4.Range('H25').Value = 'https://show.bilibili.com/platform/detail.html?id=81173'
$ws4.Range('I25').Value = '//i0.hdslb.com/bfs/openplatform/202401/TYA5FLkE1705891815532.jpeg'
# Row 28
$ws4.Range('B28').Value = '2024.03.16'
$ws4.Range('C28').Value = '上海·Look Look动漫嘉年华'
$ws4.Range('D28').Value = '龙吴路4800号2号门 有只怪兽片场'
$ws4.Range('E28').Value = '2024.03.16 10:00-03.17 17:30'
$ws4.Range('F28').Value = 66
$ws4.Range('G28').Value = 29.9
$ws4.Range('H28').Value = 'https://show.bilibili.com/platform/detail.html?id=81804'
$ws4.Range('I28').Value = '//i2.hdslb.com/bfs/openplatform/202402/WFRql6sg1707274094000.jpeg'
# Row 29
$ws4.Range('B29').Value = '2024.03.16'
$ws4.Range('C29').Value = '上海·SISP动漫游戏嘉年华'
$ws4.Range('D29').Value = '年家浜路518号 周浦万达广场'
$ws4.Range('E29').Value = '2024.03.16 13:00-03.17 19:00'
$ws4.Range('F29').Value = 167
$ws4.Range('G29').Value = 48
$ws4.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=80339'
$ws4.Range('I29').Value = '//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg'
# Row 30
$ws4.Range('B30').Value = '2024.03.16'
$ws4.Range('C30').Value = '上海·坏孩纸物语の第33届动漫节之庄子篇'
$ws4.Range('D30').Value = '中山北路3300号4楼L4001号 环球港上海世嘉都市乐园'
$ws4.Range('E30').Value = '2024.03.16 10:00-03.17 21:00'
$ws4.Range('F30').Value = 84
$ws4.Range('G30').Value = 40
$ws4.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=81138'
$ws4.Range('I30').Value = '//i2.hdslb.com/bfs/openplatform/202401/jpr1lCt21705652306481.png'
# Row 31
$ws4.Range('B31').Value = '2024.03.16'
$ws4.Range('C31').Value = '上海·第五人格ONLY'
$ws4.Range('D31').Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws4.Range('E31').Value = '2024.03.16 10:00-03.16 17:00'
$ws4.Range('F31').Value = 480
$ws4.Range('G31').Value = 60
$ws4.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=81533'
$ws4.Range('I31').Value = '//i1.hdslb.com/bfs/openplatform/202401/sOMO7Bjc1706604737277.png'
# Row 32
$ws4.Range('B32').Value = '2024.03.17'
$ws4.Range('C32').Value = '上海 ·《疯狂动物城》动漫视听音乐会'
$ws4.Range('D32').Value = '牛庄路704号 中国大戏院'
$ws4.Range('E32').Value = '2024.03.17 15:30-03.17 17:00'
$ws4.Range('F32').Value = 15
$ws4.Range('G32').Value = 80
$ws4.Range('H32').Value = 'https://show.bilibili.com/platform/detail.html?id=81112'
$ws4.Range('I32').Value = '//i2.hdslb.com/bfs/openplatform/202401/Wg8b6SRn1705651166088.png'
# Row 33
$ws4.Range('B33').Value = '2024.03.17'
$ws4.Range('C33').Value = '上海·《笑傲江湖》经典武侠影视金曲音乐会'
$ws4.Range('D33').Value = '牛庄路704号 中国大戏院'
$ws4.Range('E33').Value = '2024.03.17 19:30-03.17 21:00'
$ws4.Range('F33').Value = 3
$ws4.Range('G33').Value = 80
$ws4.Range('H33').Value = 'https://show.bilibili.com/platform/detail.html?id=80875'
$ws4.Range('I33').Value = '//i1.hdslb.com/bfs/openplatform/202401/8AwIAy4I1705385447242.jpeg'
# Row 34
$ws4.Range('B34').Value = '2024.03.17'
$ws4.Range('C34').Value = '上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会'
$ws4.Range('D34').Value = '延安东路523号 凯迪拉克·上海音乐厅'
$ws4.Range('E34').Value = '2024.03.17 14:00-03.17 16:00'
$ws4.Range('F34').Value = 53
$ws4.Range('G34').Value = 80
$ws4.Range('H34').Value = 'https://show.bilibili.com/platform/detail.html?id=81258'
$ws4.Range('I34').Value = '//i1.hdslb.com/bfs/openplatform/202401/eysvN81k1705977896972.jpeg'
# Row 35
$ws4.Range('B35').Value = '2024.03.23'
$ws4.Range('C35').Value = '上海·《卡农Canon in D》世界经典作品视听音乐会'
$ws4.Range('D35').Value = '南京西路1376号 上海商城剧院'
$ws4.Range('E35').Value = '2024.03.23 19:30-03.23 21:00'
$ws4.Range('F35').Value = 4
$ws4.Range('G35').Value = 50
$ws4.Range('H35').Value = 'https://show.bilibili.com/platform/detail.html?id=81358'
$ws4.Range('I35').Value = '//i1.hdslb.com/bfs/openplatform/202401/Ctne29Xn1706089385959.png'
# Row 36
$ws4.Range('B36').Value = '2024.03.23'
$ws4.Range('C36').Value = '上海·《四月是你的谎言》友人A经典动漫音乐会'
$ws4.Range('D36').Value = '南京西路1376号 上海商城剧院'
$ws4.Range('E36').Value = '2024.03.23 15:00-03.23 16:30'
$ws4.Range('F36').Value = 55
$ws4.Range('G36').Value = 50
$ws4.Range('H36').Value = 'https://show.bilibili.com/platform/detail.html?id=81361'
$ws4.Range('I36').Value = '//i0.hdslb.com/bfs/openplatform/202401/wL0ZWVYi1706091574963.png'
# Row 37
$ws4.Range('B37').Value = '2024.03.29'
$ws4.Range('C37').Value = '上海·KANAKO ITO&AYANE 2024 LIVE'
$ws4.Range('D37').Value = '宜昌路179号 万代南梦宫上海文化中心'
$ws4.Range('E37').Value = '2024.03.29 19:00-03.29 20:30'
$ws4.Range('F37').Value = 341
$ws4.Range('G37').Value = 380
$ws4.Range('H37').Value = 'https://show.bilibili.com/platform/detail.html?id=81416'
$ws4.Range('I37').Value = '//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg'
# Row 38
$ws4.Range('B38').Value = '2024.03.30'
$ws4.Range('C38').Value = '上海· TRUE（唐沢美帆）上海动漫交响音乐会'
$ws4.Range('D38').Value = '丁香路425号 上海东方艺术中心'
$ws4.Range('E38').Value = '2024.03.30 19:30-03.30 21:00'
$ws4.Range('F38').Value = 295
$ws4.Range('G38').Value = 680
$ws4.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=80906'
$ws4.Range('I38').Value = '//i0.hdslb.com/bfs/openplatform/202401/FaJbLvS51705401178235.jpeg'
# Row 39
$ws4.Range('B39').Value = '2024.03.30'
$ws4.Range('C39').Value = '上海·首届Redamancy动漫游戏嘉年华'
$ws4.Range('D39').Value = '中山北路3300号4楼L4001号 环球港上海世嘉都市乐园'
$ws4.Range('E39').Value = '2024.03.30 10:00-03.31 17:00'
$ws4.Range('F39').Value = 1262
$ws4.Range('G39').Value = 60
$ws4.Range('H39').Value = 'https://show.bilibili.com/platform/detail.html?id=81772'
$ws4.Range('I39').Value = '//i2.hdslb.com/bfs/openplatform/202402/XKf9RSFB1707127784856.jpeg'
# Row 40
$ws4.Range('B40').Value = '2024.04.06'
$ws4.Range('C40').Value = '上海·从Butter-Fly到夏目之爱してる —— “好想大声说爱你”动漫钢琴演奏会'
$ws4.Range('D40').Value = '复兴中路1380号 捷豹上海交响音乐厅'
$ws4.Range('E40').Value = '2024.04.06 19:30-04.06 21:30'
$ws4.Range('F40').Value = 25
$ws4.Range('G40').Value = 80
$ws4.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=80050'
$ws4.Range('I40').Value = '//i0.hdslb.com/bfs/openplatform/202312/0iJP3TY61703056498448.jpeg'
# Row 41
$ws4.Range('B41').Value = '2024.04.13'
$ws4.Range('C41').Value = '上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集'
$ws4.Range('D41').Value = '丁香路425号 上海东方艺术中心'
$ws4.Range('E41').Value = '2024.04.13 19:30-04.13 21:30'
$ws4.Range('F41').Value = 210
$ws4.Range('G41').Value = 80
$ws4.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=78667'
$ws4.Range('I41').Value = '//i1.hdslb.com/bfs/openplatform/202311/bTP7w6GD1700130122940.jpeg'
# Row 42
$ws4.Range('B42').Value = '2024.04.13'
$ws4.Range('C42').Value = '上海·运动番ONLY'
$ws4.Range('D42').Value = '少年村路6号 YC篮羽联盟(大场店)'
$ws4.Range('E42').Value = '2024.04.13 10:00-04.13 17:00'
$ws4.Range('F42').Value = 10
$ws4.Range('G42').Value = 60
$ws4.Range('H42').Value = 'https://show.bilibili.com/platform/detail.html?id=81901'
$ws4.Range('I42').Value = '//i0.hdslb.com/bfs/openplatform/202402/2oiNlCAr1708325440584.jpeg'

# --- Sheet 4: delete row 50 (event dropped); dimension becomes A1:I49 ---
$ws4.Rows.Item(50).Delete()
